# Scheduled runner update: refresh computed profit-margin figures across
# all item sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with the latest
# market-snapshot numbers (columns H..N per changed row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2898
$ws.Range("I32").Value = 2499
$ws.Range("J32").Value = 2997.75
$ws.Range("K32").Value = 2499
$ws.Range("L32").Value = 2997.75
$ws.Range("M32").Value = -2173
$ws.Range("N32").Value = -3649.75

$ws.Range("H38").Value = 401.57144
$ws.Range("I38").Value = 340.15384
$ws.Range("J38").Value = 1200
$ws.Range("K38").Value = 1020.46152
$ws.Range("L38").Value = 3600
$ws.Range("M38").Value = -648.4615200000001
$ws.Range("N38").Value = -4344

$ws.Range("H86").Value = 2036.0769
$ws.Range("I86").Value = 1809.375
$ws.Range("K86").Value = 1809.375
$ws.Range("M86").Value = -686.375

$ws.Range("H89").Value = 2036.0769
$ws.Range("I89").Value = 1809.375
$ws.Range("K89").Value = 9046.875
$ws.Range("M89").Value = -3430.875

$ws.Range("H98").Value = 1796.7333
$ws.Range("I98").Value = 1079.5
$ws.Range("K98").Value = 1079.5
$ws.Range("M98").Value = 418.5

$ws.Range("H122").Value = 1796.7333
$ws.Range("I122").Value = 1079.5
$ws.Range("K122").Value = 3238.5
$ws.Range("M122").Value = -788.5

$ws.Range("H125").Value = 1380.8235
$ws.Range("J125").Value = 1826
$ws.Range("L125").Value = 16434
$ws.Range("N125").Value = -21354

$ws.Range("H135").Value = 813.55554
$ws.Range("I135").Value = 813.55554
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 7321.99986
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -4786.99986
$ws.Range("N135").ClearContents()

$ws.Range("H138").Value = 2857.04
$ws.Range("I138").Value = 2094.3
$ws.Range("K138").Value = 6282.900000000001
$ws.Range("M138").Value = -1142.900000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2676.2632
$ws.Range("I2").Value = 2285.8823
$ws.Range("K2").Value = 2285.8823
$ws.Range("M2").Value = -2172.8823

$ws.Range("H4").Value = 200707.9
$ws.Range("I4").Value = 250419.88
$ws.Range("J4").Value = 1860
$ws.Range("K4").Value = 250419.88
$ws.Range("L4").Value = 1860
$ws.Range("M4").Value = -250303.88
$ws.Range("N4").Value = -2092

$ws.Range("H32").Value = 3214.5066
$ws.Range("I32").Value = 2993.6448
$ws.Range("K32").Value = 2993.6448
$ws.Range("M32").Value = -2706.6448

$ws.Range("H50").Value = 748
$ws.Range("I50").Value = 597.5
$ws.Range("J50").Value = 1350
$ws.Range("K50").Value = 597.5
$ws.Range("L50").Value = 1350
$ws.Range("M50").Value = 116.5
$ws.Range("N50").Value = -2778

$ws.Range("H57").Value = 4975
$ws.Range("I57").Value = 4975
$ws.Range("K57").Value = 4975
$ws.Range("M57").Value = -4491

$ws.Range("H63").Value = 3721.111
$ws.Range("I63").Value = 2936.25
$ws.Range("K63").Value = 2936.25
$ws.Range("M63").Value = -2250.25

$ws.Range("H66").Value = 3721.111
$ws.Range("I66").Value = 2936.25
$ws.Range("K66").Value = 14681.25
$ws.Range("M66").Value = -11249.25

$ws.Range("H97").Value = 4281.905
$ws.Range("I97").Value = 2026.8422
$ws.Range("J97").Value = 25705
$ws.Range("K97").Value = 2026.8422
$ws.Range("L97").Value = 25705
$ws.Range("M97").Value = -1530.8422
$ws.Range("N97").Value = -26697

$ws.Range("H116").Value = 2676.2632
$ws.Range("I116").Value = 2285.8823
$ws.Range("K116").Value = 2285.8823
$ws.Range("M116").Value = 8.117699999999786

$ws.Range("H126").Value = 5010.75
$ws.Range("I126").Value = 5010.75
$ws.Range("K126").Value = 15032.25
$ws.Range("M126").Value = -12562.25

$ws.Range("H132").Value = 6600.162
$ws.Range("I132").Value = 6600.162
$ws.Range("K132").Value = 19800.486
$ws.Range("M132").Value = -17270.486

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2676.2632
$ws.Range("I3").Value = 2285.8823
$ws.Range("K3").Value = 2285.8823
$ws.Range("M3").Value = -2171.8823

$ws.Range("H20").Value = 6010.1387
$ws.Range("I20").Value = 7527.231
$ws.Range("J20").Value = 2065.7
$ws.Range("K20").Value = 7527.231
$ws.Range("L20").Value = 2065.7
$ws.Range("M20").Value = -7280.231
$ws.Range("N20").Value = -2559.7

$ws.Range("H128").Value = 19646.834
$ws.Range("I128").Value = 19646.834
$ws.Range("K128").Value = 58940.50199999999
$ws.Range("M128").Value = -56450.50199999999

$ws.Range("H134").Value = 167967.17
$ws.Range("I134").Value = 167967.17
$ws.Range("K134").Value = 503901.51
$ws.Range("M134").Value = -501366.51

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3960.5938
$ws.Range("J31").Value = 7312.6665
$ws.Range("L31").Value = 7312.6665
$ws.Range("N31").Value = -7902.6665

$ws.Range("H34").Value = 3960.5938
$ws.Range("J34").Value = 7312.6665
$ws.Range("L34").Value = 7312.6665
$ws.Range("N34").Value = -7716.6665

$ws.Range("H86").Value = 3086.25
$ws.Range("I86").Value = 2943.6924
$ws.Range("K86").Value = 2943.6924
$ws.Range("M86").Value = -1820.6924

$ws.Range("H89").Value = 3086.25
$ws.Range("I89").Value = 2943.6924
$ws.Range("K89").Value = 14718.462
$ws.Range("M89").Value = -9102.462

$ws.Range("H107").Value = 1986.1428
$ws.Range("I107").Value = 1626.7059
$ws.Range("K107").Value = 1626.7059
$ws.Range("M107").Value = 293.2941000000001

$ws.Range("H132").Value = 3529.6428
$ws.Range("I132").Value = 3562.6924
$ws.Range("J132").Value = 3100
$ws.Range("K132").Value = 10688.0772
$ws.Range("L132").Value = 9300
$ws.Range("M132").Value = -8158.0772
$ws.Range("N132").Value = -14360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2301890.2
$ws.Range("I4").Value = 1264212.1
$ws.Range("J4").Value = 3487808
$ws.Range("K4").Value = 3792636.3
$ws.Range("L4").Value = 10463424
$ws.Range("M4").Value = -3792524.3
$ws.Range("N4").Value = -10463648

$ws.Range("H113").Value = 2328.85
$ws.Range("I113").Value = 701
$ws.Range("J113").Value = 2616.1177
$ws.Range("K113").Value = 2103
$ws.Range("L113").Value = 7848.353099999999
$ws.Range("M113").Value = 67
$ws.Range("N113").Value = -12188.3531

$ws.Range("H137").Value = 3821.5
$ws.Range("J137").Value = 4219.0625
$ws.Range("L137").Value = 12657.1875
$ws.Range("N137").Value = -22857.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3187.4666
$ws.Range("J80").Value = 4080.4285
$ws.Range("L80").Value = 4080.4285
$ws.Range("N80").Value = -6076.4285

$ws.Range("H83").Value = 3187.4666
$ws.Range("J83").Value = 4080.4285
$ws.Range("L83").Value = 20402.1425
$ws.Range("N83").Value = -30386.1425

$ws.Range("H113").Value = 4747.5
$ws.Range("I113").Value = 3500
$ws.Range("K113").Value = 3500
$ws.Range("M113").Value = -1330

$ws.Range("H122").Value = 2783.8948
$ws.Range("I122").Value = 2549.6428
$ws.Range("J122").Value = 3439.8
$ws.Range("K122").Value = 7648.928400000001
$ws.Range("L122").Value = 10319.4
$ws.Range("M122").Value = -5198.928400000001
$ws.Range("N122").Value = -15219.4

$ws.Range("H123").Value = 66738.91
$ws.Range("J123").Value = 66738.91
$ws.Range("L123").Value = 66738.91
$ws.Range("N123").Value = -71638.91

$ws.Range("H126").Value = 2342.889
$ws.Range("I126").Value = 2324.7693
$ws.Range("K126").Value = 6974.3079
$ws.Range("M126").Value = -4504.3079

$ws.Range("H132").Value = 3419.849
$ws.Range("I132").Value = 2616
$ws.Range("J132").Value = 6876.4
$ws.Range("K132").Value = 7848
$ws.Range("L132").Value = 20629.2
$ws.Range("M132").Value = -5318
$ws.Range("N132").Value = -25689.2

$ws.Range("H136").Value = 60417.05
$ws.Range("J136").Value = 60417.05
$ws.Range("L136").Value = 181251.15
$ws.Range("N136").Value = -186351.15

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2961.524
$ws.Range("I136").Value = 2612.2942
$ws.Range("J136").Value = 4445.75
$ws.Range("K136").Value = 7836.882599999999
$ws.Range("L136").Value = 13337.25
$ws.Range("M136").Value = -5286.882599999999
$ws.Range("N136").Value = -18437.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11149.533
$ws.Range("I62").Value = 4353.5835
$ws.Range("K62").Value = 4353.5835
$ws.Range("M62").Value = -3729.5835

$ws.Range("H65").Value = 11149.533
$ws.Range("I65").Value = 4353.5835
$ws.Range("K65").Value = 21767.9175
$ws.Range("M65").Value = -18647.9175

$ws.Range("H107").Value = 62513516
$ws.Range("I107").Value = 15448.857
$ws.Range("J107").Value = 500000000
$ws.Range("K107").Value = 46346.571
$ws.Range("L107").Value = 1500000000
$ws.Range("M107").Value = -44426.571
$ws.Range("N107").Value = -1500003840

$ws.Range("H122").Value = 1529.6136
$ws.Range("I122").Value = 1508.9375
$ws.Range("J122").Value = 1584.75
$ws.Range("K122").Value = 4526.8125
$ws.Range("L122").Value = 4754.25
$ws.Range("M122").Value = -2076.8125
$ws.Range("N122").Value = -9654.25

$ws.Range("H132").Value = 1164772.9
$ws.Range("I132").Value = 2955.96
$ws.Range("K132").Value = 8867.880000000001
$ws.Range("M132").Value = -6337.880000000001
